$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 92
$ws.Cells.Item(2, 6).Value = 64
$ws.Cells.Item(2, 8).Value = 64

$ws.Cells.Item(5, 5).Value = 132

$ws.Cells.Item(10, 5).Value = 490
$ws.Cells.Item(10, 6).Value = 248
$ws.Cells.Item(10, 8).Value = 248

$ws.Cells.Item(11, 5).Value = 325
$ws.Cells.Item(11, 6).Value = 185
$ws.Cells.Item(11, 8).Value = 185

$ws.Cells.Item(12, 5).Value = 478
$ws.Cells.Item(12, 6).Value = 261
$ws.Cells.Item(12, 8).Value = 261

$ws.Cells.Item(13, 5).Value = 124
$ws.Cells.Item(13, 6).Value = 67
$ws.Cells.Item(13, 8).Value = 67

$ws.Cells.Item(14, 5).Value = 121
$ws.Cells.Item(14, 6).Value = 64
$ws.Cells.Item(14, 8).Value = 64

$ws.Cells.Item(15, 5).Value = 159
$ws.Cells.Item(15, 6).Value = 68
$ws.Cells.Item(15, 8).Value = 68

$ws.Cells.Item(16, 5).Value = 192

$ws.Cells.Item(17, 5).Value = 94
$ws.Cells.Item(17, 6).Value = 49
$ws.Cells.Item(17, 8).Value = 49

$ws.Cells.Item(21, 5).Value = 133
$ws.Cells.Item(21, 6).Value = 73
$ws.Cells.Item(21, 8).Value = 73

$ws.Cells.Item(23, 5).Value = 194
$ws.Cells.Item(23, 6).Value = 91
$ws.Cells.Item(23, 8).Value = 91

$ws.Cells.Item(24, 5).Value = 203

$ws.Cells.Item(27, 5).Value = 313

$ws.Cells.Item(29, 5).Value = 164

$ws.Cells.Item(30, 5).Value = 195
$ws.Cells.Item(30, 6).Value = 117
$ws.Cells.Item(30, 8).Value = 117

$ws.Cells.Item(31, 5).Value = 71
$ws.Cells.Item(31, 6).Value = 32
$ws.Cells.Item(31, 8).Value = 32

$ws.Cells.Item(32, 5).Value = 179
$ws.Cells.Item(32, 6).Value = 106
$ws.Cells.Item(32, 8).Value = 106

$ws.Cells.Item(33, 5).Value = 279
$ws.Cells.Item(33, 6).Value = 142
$ws.Cells.Item(33, 8).Value = 142

$ws.Cells.Item(34, 6).Value = 136
$ws.Cells.Item(34, 8).Value = 136

$ws.Cells.Item(36, 5).Value = 67

$ws.Cells.Item(39, 5).Value = 176
$ws.Cells.Item(39, 6).Value = 87
$ws.Cells.Item(39, 8).Value = 87

$ws.Cells.Item(40, 5).Value = 252
$ws.Cells.Item(40, 6).Value = 120
$ws.Cells.Item(40, 8).Value = 120

$ws.Cells.Item(42, 5).Value = 358
$ws.Cells.Item(42, 6).Value = 198
$ws.Cells.Item(42, 8).Value = 198

$ws.Cells.Item(43, 5).Value = 111
$ws.Cells.Item(43, 6).Value = 61
$ws.Cells.Item(43, 8).Value = 61

$ws.Cells.Item(45, 5).Value = 135

$ws.Cells.Item(47, 5).Value = 430
$ws.Cells.Item(47, 6).Value = 213
$ws.Cells.Item(47, 8).Value = 213

$ws.Cells.Item(48, 5).Value = 195

$ws.Cells.Item(50, 5).Value = 234

$ws.Cells.Item(51, 5).Value = 225
$ws.Cells.Item(51, 6).Value = 97
$ws.Cells.Item(51, 8).Value = 97

$ws.Cells.Item(52, 5).Value = 25
$ws.Cells.Item(52, 6).Value = 12
$ws.Cells.Item(52, 8).Value = 12

